$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row (row 1) category labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Update the distractor/target indicator values so that n distractor = n targets
$values = @{
    2 = @(0, 0, 0, 1, 0, 0)
    3 = @(0, 1, 0, 0, 0, 0)
    4 = @(0, 0, 0, 0, 0, 1)
    5 = @(0, 0, 1, 0, 0, 0)
    6 = @(1, 0, 0, 0, 0, 0)
    7 = @(0, 0, 0, 0, 1, 0)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($row, $c).Value = $rowVals[$c - 1]
    }
}
